# October 7th Work. (#72)
# Fill in the day's progress entry for row 8 (2025-10-07 -> serial 45937)
# with the recorded skills/talents/etc., matching the formatting
# (Neutral/Good/Bad cell styles) already used on the other daily rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 currently holds only the date (A8) with all other cells blank and
# styled with the generic bordered style (s=3). Copy the cell formatting
# (fill/border/number format) from the previous day's row (row 7) onto the
# same columns of row 8 so B8/F8/J8 become "Neutral", C8/E8/G8/H8/I8 become
# "Good", and D8/K8/L8/M8 become "Bad", exactly like every other data row.
$ws.Range("B7:M7").Copy()
$ws.Range("B8:M8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Now populate the actual values for the day. The order below matches the
# order the strings were originally typed into the workbook (and therefore
# the order they were appended to the shared-strings table).
$ws.Range("J8").Value = "Augure"
$ws.Range("B8").Value = "Hermit"
$ws.Range("C8").Value = "Initiative"
$ws.Range("G8").Value = "Combinaison élémentaire"
$ws.Range("H8").Value = "Endurance accrue"
$ws.Range("F8").Value = "Entreprise familiale"
$ws.Range("E8").Value = "Chaotique"
$ws.Range("I8").Value = "Attaque puissante"
